$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 112307579
$ws.Range("B5").Value = 96735
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 220787
$ws.Range("F5").Value = "Knärot"
$ws.Range("G5").Value = "Goodyera repens"
$ws.Range("H5").Value = "(L.) R. Br."
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "10"
$ws.Range("J5").Value = "plantor/tuvor"
$ws.Range("L5").NumberFormat = "General"
$ws.Range("Q5").Value = 683344
$ws.Range("R5").Value = 6627679

# Row 6
$ws.Range("A6").Value = 112307588
$ws.Range("B6").Value = 89950
$ws.Range("E6").Value = 5420
$ws.Range("F6").Value = "Grovticka"
$ws.Range("G6").Value = "Phaeolus schweinitzii"
$ws.Range("H6").Value = "(Fr.) Pat."
$ws.Range("Q6").Value = 683341
$ws.Range("R6").Value = 6627677

# Row 7
$ws.Range("A7").Value = 112307555
$ws.Range("B7").Value = 89331
$ws.Range("E7").Value = 3215
$ws.Range("F7").Value = "Rödgul trumpetsvamp"
$ws.Range("G7").Value = "Craterellus lutescens"
$ws.Range("H7").Value = "(Fr.) Fr."
$ws.Range("I7").NumberFormat = "General"
$ws.Range("J7").NumberFormat = "General"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("Q7").Value = 683437
$ws.Range("AC7").Value = ""

# Row 8
$ws.Range("A8").Value = 112307976
$ws.Range("B8").Value = 90826
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4366
$ws.Range("F8").Value = "Skarp dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum peckii"
$ws.Range("H8").Value = "Banker"
$ws.Range("Q8").Value = 683370
$ws.Range("R8").Value = 6627485

# Row 9
$ws.Range("A9").Value = 112307600
$ws.Range("B9").Value = 90806
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 4361
$ws.Range("F9").Value = "Orange taggsvamp"
$ws.Range("G9").Value = "Hydnellum aurantiacum"
$ws.Range("H9").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q9").Value = 683355
$ws.Range("R9").Value = 6627623

# Row 10
$ws.Range("A10").Value = 112307576
$ws.Range("B10").Value = 89553
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = "Ullticka"
$ws.Range("G10").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H10").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q10").Value = 683375
$ws.Range("R10").Value = 6627671

# Row 11
$ws.Range("A11").Value = 112307952
$ws.Range("B11").Value = 89553
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = "Ullticka"
$ws.Range("G11").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H11").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("L11").Value = ""
$ws.Range("Q11").Value = 683399
$ws.Range("R11").Value = 6627533

# Row 12
$ws.Range("A12").Value = 112307580
$ws.Range("B12").Value = 5113
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 100526
$ws.Range("F12").Value = "Bronshjon"
$ws.Range("G12").Value = "Callidium coriaceum"
$ws.Range("H12").Value = "Paykull, 1800"
$ws.Range("L12").NumberFormat = "General"
$ws.Range("M12").NumberFormat = "General"
$ws.Range("Q12").Value = 683344
$ws.Range("R12").Value = 6627679

# Row 13
$ws.Range("A13").Value = 112307943
$ws.Range("B13").Value = 89573
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 5442
$ws.Range("F13").Value = "Tallticka"
$ws.Range("G13").Value = "Porodaedalea pini"
$ws.Range("H13").Value = "(Brot.) Murrill"
$ws.Range("L13").Value = ""
$ws.Range("Q13").Value = 683401
$ws.Range("R13").Value = 6627535

# Row 14
$ws.Range("A14").Value = 112307990
$ws.Range("I14").NumberFormat = "General"
$ws.Range("J14").NumberFormat = "General"
$ws.Range("Q14").Value = 683334
$ws.Range("R14").Value = 6627502

# Row 15
$ws.Range("A15").Value = 112307609
$ws.Range("B15").Value = 90822
$ws.Range("D15").Value = "VU"
$ws.Range("E15").Value = 2058
$ws.Range("F15").Value = "Koppartaggsvamp"
$ws.Range("G15").Value = "Hydnellum lundellii"
$ws.Range("H15").Value = "(Maas Geest. & Nannf.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q15").Value = 683391
$ws.Range("R15").Value = 6627583

# Row 16
$ws.Range("A16").Value = 112307568
$ws.Range("B16").Value = 96735
$ws.Range("D16").Value = "VU"
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = "Knärot"
$ws.Range("G16").Value = "Goodyera repens"
$ws.Range("H16").Value = "(L.) R. Br."
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = "10"
$ws.Range("J16").Value = "plantor/tuvor"
$ws.Range("L16").NumberFormat = "General"
$ws.Range("Q16").Value = 683342
$ws.Range("R16").Value = 6627694

# Row 17
$ws.Range("A17").Value = 112307572
$ws.Range("B17").Value = 89993
$ws.Range("E17").Value = 1209
$ws.Range("F17").Value = "Rynkskinn"
$ws.Range("G17").Value = "Phlebia centrifuga"
$ws.Range("H17").Value = "P.Karst."
$ws.Range("I17").NumberFormat = "General"
$ws.Range("J17").NumberFormat = "General"
$ws.Range("L17").Value = ""
$ws.Range("Q17").Value = 683375
$ws.Range("R17").Value = 6627671

# Row 18
$ws.Range("A18").Value = 112307970
$ws.Range("B18").Value = 90806
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 4361
$ws.Range("F18").Value = "Orange taggsvamp"
$ws.Range("G18").Value = "Hydnellum aurantiacum"
$ws.Range("H18").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q18").Value = 683369
$ws.Range("R18").Value = 6627485

# Row 19
$ws.Range("A19").Value = 112307534
$ws.Range("B19").Value = 90837
$ws.Range("E19").Value = 5966
$ws.Range("F19").Value = "Motaggsvamp"
$ws.Range("G19").Value = "Sarcodon squamosus"
$ws.Range("H19").Value = "(Schaeff.) Quél."
$ws.Range("Q19").Value = 683417
$ws.Range("R19").Value = 6627694

# Row 20
$ws.Range("A20").Value = 112307592
$ws.Range("B20").Value = 103781
$ws.Range("E20").Value = 221144
$ws.Range("F20").Value = "Grönpyrola"
$ws.Range("G20").Value = "Pyrola chlorantha"
$ws.Range("H20").Value = "Sw."
$ws.Range("L20").NumberFormat = "General"
$ws.Range("Q20").Value = 683347
$ws.Range("R20").Value = 6627644

# Row 22
$ws.Range("A22").Value = 112307503
$ws.Range("B22").Value = 90806
$ws.Range("E22").Value = 4361
$ws.Range("F22").Value = "Orange taggsvamp"
$ws.Range("G22").Value = "Hydnellum aurantiacum"
$ws.Range("H22").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q22").Value = 683373
$ws.Range("R22").Value = 6627724

# Row 23
$ws.Range("A23").Value = 112307522
$ws.Range("B23").Value = 90826
$ws.Range("D23").Value = "LC"
$ws.Range("E23").Value = 4366
$ws.Range("F23").Value = "Skarp dropptaggsvamp"
$ws.Range("G23").Value = "Hydnellum peckii"
$ws.Range("H23").Value = "Banker"
$ws.Range("Q23").Value = 683384
$ws.Range("R23").Value = 6627716
